$d = $word.ActiveDocument

# Locate the two target paragraphs by content instead of a hard-coded
# index, so the script is robust to the exact paragraph numbering.
$n = $d.Paragraphs.Count
$idxIntro = -1
$idxShould = -1
for ($i = 1; $i -le $n; $i++) {
  $txt = $d.Paragraphs.Item($i).Range.Text
  if ($txt -like "*I am in the process of preparing a manuscript*") {
    $idxIntro = $i
  }
  if ($txt -like "*Should you choose to review and contribute*") {
    $idxShould = $i
  }
}

# --- Change 1 & 2 ---------------------------------------------------
# The "I am in the process ... different groups." paragraph:
#   - pPr/rPr rFonts hint: eastAsia -> default
#   - run text gets a new trailing sentence appended
$pIntro = $d.Paragraphs.Item($idxIntro)
$xmlIntro = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" + `
  "<w:pPr><w:rPr><w:rFonts w:hint='default'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr></w:pPr>" + `
  "<w:r><w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr>" + `
  "<w:t>I am in the process of preparing a manuscript for submission to PNAS. The paper presents three novel types of dissimilarity measure that leverage moment and sparsity differences. This measure can be applied across various omics datasets, offering a powerful tool for uncovering the pivotal factors that distinguish different groups. I am particularly eager to bring to your attention this paper that I believe will resonate with your field of study.</w:t>" + `
  "</w:r></w:p>"
$pIntro.Range.InsertXML($xmlIntro)

# --- Change 3 ---------------------------------------------------------
# The "Should you choose to review..." paragraph plus the following
# (empty, bookmark-carrying) paragraph get restructured: a new sentence
# is inserted as its own run before the "Should you..." run, the
# "_GoBack" bookmark moves in between the two runs, and the trailing
# empty paragraph loses the bookmark (keeping only its pPr).
$pShould = $d.Paragraphs.Item($idxShould)
$pAfter = $d.Paragraphs.Item($idxShould + 1)
$start = $pShould.Range.Start
$end = $pAfter.Range.End
$r = $d.Range($start, $end)

$rPr = "<w:rPr><w:rFonts w:hint='eastAsia'/><w:lang w:val='en-US' w:eastAsia='zh-CN'/></w:rPr>"
$pPr = "<w:pPr>$rPr</w:pPr>"

$xmlShould = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>$pPr" + `
  "<w:r>$rPr<w:t xml:space='preserve'>In support of our collaborative efforts, I have attached both the source code and manuscripts. I encourage you to apply these methods to your previously published datasets to explore their utility. </w:t></w:r>" + `
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" + `
  "<w:r>$rPr<w:t>Should you choose to review and contribute to that paper, I would be honored to include you as an author.</w:t></w:r>" + `
  "</w:p>" + `
  "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>$pPr</w:p>"

$r.InsertXML($xmlShould)

Write-Output "edit applied"
